$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRange, $value) {
    $origStyle = $cellRange.Style
    $cellRange.NumberFormat = "@"
    $cellRange.Value = $value
    $cellRange.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "30.673.26"
Set-TextValue $ws.Range("E2") "  +1.88%  "

Set-TextValue $ws.Range("D3") "1.893.48"
Set-TextValue $ws.Range("E3") "  +1.06%  "

Set-TextValue $ws.Range("D4") "1.001"
Set-TextValue $ws.Range("E4") "  +0.06%  "

Set-TextValue $ws.Range("D5") "241.61"
Set-TextValue $ws.Range("E5") "  -0.40%  "

Set-TextValue $ws.Range("E6") "  +0.11%  "

Set-TextValue $ws.Range("D7") "0.4908"
Set-TextValue $ws.Range("E7") "  +0.69%  "

Set-TextValue $ws.Range("E8") "  +1.47%  "

Set-TextValue $ws.Range("D9") "0.06738"
Set-TextValue $ws.Range("E9") "  +2.53%  "

Set-TextValue $ws.Range("D10") "1.893.96"
Set-TextValue $ws.Range("E10") "  +1.13%  "

Set-TextValue $ws.Range("D11") "17.13"
Set-TextValue $ws.Range("E11") "  +5.23%  "

Set-TextValue $ws.Range("D12") "0.07259"
Set-TextValue $ws.Range("E12") "  +0.82%  "

Set-TextValue $ws.Range("D13") "90.74"
Set-TextValue $ws.Range("E13") "  +5.77%  "

Set-TextValue $ws.Range("D14") "0.6745"
Set-TextValue $ws.Range("E14") "  +1.87%  "

Set-TextValue $ws.Range("D15") "5.019"
Set-TextValue $ws.Range("E15") "  +1.57%  "

Set-TextValue $ws.Range("D16") "30.659.84"
Set-TextValue $ws.Range("E16") "  +2.03%  "

Set-TextValue $ws.Range("D17") "0.000007954"
Set-TextValue $ws.Range("E17") "  +2.61%  "

Set-TextValue $ws.Range("E18") "  +0.18%  "

Set-TextValue $ws.Range("D19") "13.10"
Set-TextValue $ws.Range("E19") "  +2.85%  "

Set-TextValue $ws.Range("D20") "2.139.42"
Set-TextValue $ws.Range("E20") "  +1.20%  "

Set-TextValue $ws.Range("D21") "1.003"
Set-TextValue $ws.Range("E21") "  +0.52%  "

Set-TextValue $ws.Range("E22") "  +1.05%  "

Set-TextValue $ws.Range("D23") "191.71"
Set-TextValue $ws.Range("E23") "  +33.62%  "

Set-TextValue $ws.Range("D24") "6.075"
Set-TextValue $ws.Range("E24") "  +3.51%  "

Set-TextValue $ws.Range("D25") "9.358"
Set-TextValue $ws.Range("E25") "  +2.33%  "

Set-TextValue $ws.Range("D26") "157.11"
Set-TextValue $ws.Range("E26") "  +3.69%  "

Set-TextValue $ws.Range("D27") "18.85"
Set-TextValue $ws.Range("E27") "  +11.52%  "

Set-TextValue $ws.Range("D28") "1.893"
Set-TextValue $ws.Range("E28") "  +0.81%  "

Set-TextValue $ws.Range("D29") "1.403"
Set-TextValue $ws.Range("E29") "  +0.54%  "

Set-TextValue $ws.Range("D30") "4.288"
Set-TextValue $ws.Range("E30") "  +2.39%  "

Set-TextValue $ws.Range("D31") "0.09068"
Set-TextValue $ws.Range("E31") "  +3.38%  "

Set-TextValue $ws.Range("D32") "3.992"

Set-TextValue $ws.Range("D33") "0.05238"
Set-TextValue $ws.Range("E33") "  +1.27%  "

Set-TextValue $ws.Range("D34") "0.7372"
Set-TextValue $ws.Range("E34") "  +2.89%  "

Set-TextValue $ws.Range("E35") "  +0.00%  "

Set-TextValue $ws.Range("D36") "2.727"
Set-TextValue $ws.Range("E36") "  +2.37%  "

Set-TextValue $ws.Range("E37") "  -0.50%  "

Set-TextValue $ws.Range("D38") "2.677"
Set-TextValue $ws.Range("E38") "  +0.28%  "

Set-TextValue $ws.Range("D39") "0.9315"
Set-TextValue $ws.Range("E39") "  +0.54%  "

Set-TextValue $ws.Range("D40") "2.117"
Set-TextValue $ws.Range("E40") "  -1.80%  "

Set-TextValue $ws.Range("D41") "0.4391"
Set-TextValue $ws.Range("E41") "  +3.61%  "

Set-TextValue $ws.Range("D42") "104.94"
Set-TextValue $ws.Range("E42") "  +1.28%  "

Set-TextValue $ws.Range("E43") "  +0.22%  "

Set-TextValue $ws.Range("D44") "5.728"
Set-TextValue $ws.Range("E44") "  -0.33%  "

Set-TextValue $ws.Range("E45") "  +1.15%  "

Set-TextValue $ws.Range("D46") "0.1347"
Set-TextValue $ws.Range("E46") "  +5.29%  "

Set-TextValue $ws.Range("D47") "0.05862"
Set-TextValue $ws.Range("E47") "  +2.47%  "

Set-TextValue $ws.Range("B48") "NEARProtocol"
Set-TextValue $ws.Range("C48") "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Range("D48") "1.431"
Set-TextValue $ws.Range("E48") "  +6.82%  "

Set-TextValue $ws.Range("B49") "EnergySwap"
Set-TextValue $ws.Range("C49") "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D49") "8.689"
Set-TextValue $ws.Range("E49") "  +5.30%  "

Set-TextValue $ws.Range("D50") "0.3937"
Set-TextValue $ws.Range("E50") "  +4.87%  "

Set-TextValue $ws.Range("D51") "33.69"
Set-TextValue $ws.Range("E51") "  +2.86%  "
